# Weekly update: insert a new daily record for "Brócoli" at
# Macroferia Regional de Talca, pushing the existing rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 232 (row 1 is the header; data starts at row 2).
# Everything that was in row 232..294 shifts down to 233..295.
$ws.Rows("232:232").Insert()

# Populate the newly inserted row with the new record's values.
$ws.Range("A232").Value2 = 5
$ws.Range("B232").Value2 = "Macroferia Regional de Talca"
$ws.Range("C232").Value2 = "Maule"
$ws.Range("D232").Value2 = 44642
$ws.Range("E232").Value2 = 7
$ws.Range("F232").Value2 = 100112023
$ws.Range("G232").Value2 = "Brócoli"
$ws.Range("H232").Value2 = "Sin especificar"
$ws.Range("I232").Value2 = "Primera"
$ws.Range("J232").Value2 = 4000
$ws.Range("K232").Value2 = 700
$ws.Range("L232").Value2 = 700
$ws.Range("M232").Value2 = 700
$ws.Range("N232").Value2 = '$/unidad'
$ws.Range("O232").Value2 = "Región del Maule"
$ws.Range("P232").Value2 = 700
$ws.Range("Q232").Value2 = 1
$ws.Range("R232").Value2 = "Hortaliza"

# Match the date cell's number format used by the rest of column D.
$ws.Range("D232").NumberFormat = $ws.Range("D233").NumberFormat
